$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column values for rows 2-7 from 45175 to 45183.
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 3).Value = 45183
}
